# Slide 3's title shape ("Title 1") currently holds its text as three
# separate runs: "Below", " ", "section-level" (all sharing the same,
# empty run properties). Merge them into a single run "Below section-level"
# while leaving the run properties untouched (still an empty <a:rPr/>).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange

# Append the remaining words onto the end of the first run ("Below") so
# the extra text is absorbed into that run rather than creating a new
# one with freshly generated run properties.
$firstRun = $tr.Characters(1, 5)
[void]$firstRun.InsertAfter(" section-level")

# The original second/third runs' text still follows; remove that
# now-duplicated trailing text.
$tail = $tr.Characters(20, 14)
$tail.Text = ""
